$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title block: merge B2:H3 and add "Class Diagram" title with dark-blue fill ---
$ws.Range("B2:H3").Merge()
$ws.Range("B2").Value = "Class Diagram"
$ws.Range("B2:H3").Interior.Color = 6299648

# --- New "Database" class block (rows 22-31), styled like the other class blocks ---
$ws.Range("B5:B14").Copy()
$ws.Range("B22").PasteSpecial(-4122)

$ws.Range("B22").Value = "Database"
$ws.Range("B23").Value = "[ - ] host : String"
$ws.Range("B24").Value = "[ - ] user : String"
$ws.Range("B25").Value = "[ - ] passw : String"
$ws.Range("B26").Value = "[ - ] db : String"
$ws.Range("B27").Value = "[ - ] conn : String"
$ws.Range("B29").Value = "[ - ] linkage() : void"
$ws.Range("B30").Value = "[ # ] connect() : void"

$ws.Rows(22).RowHeight = 15.5
$ws.Rows(28).RowHeight = 15

# --- Update view state to match the author's last selection/scroll position ---
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("C30").Select()
